# Scheduled-runner refresh: overwrite cached FFXIV market-price / leve-profit
# figures (columns H-N) with the latest pulled values. No formulas are
# involved anywhere in this workbook -- every cell below is a plain number
# literal, so this script just re-pokes the updated values cell-by-cell,
# sheet-by-sheet (matching the commit's per-row diffs).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 58898684  # H76
$ws.Cells.Item(76, 10).Value = 90927560  # J76
$ws.Cells.Item(76, 12).Value = 90927560  # L76
$ws.Cells.Item(76, 14).Value = -90928190  # N76

$ws.Cells.Item(79, 8).Value = 58898684  # H79
$ws.Cells.Item(79, 10).Value = 90927560  # J79
$ws.Cells.Item(79, 12).Value = 90927560  # L79
$ws.Cells.Item(79, 14).Value = -90929744  # N79

$ws.Cells.Item(138, 8).Value = 3779.9666  # H138
$ws.Cells.Item(138, 10).Value = 5648.5  # J138
$ws.Cells.Item(138, 12).Value = 16945.5  # L138
$ws.Cells.Item(138, 14).Value = -27225.5  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 100000500  # H4
$ws.Cells.Item(4, 9).Value = 424.85715  # I4
$ws.Cells.Item(4, 11).Value = 424.85715  # K4
$ws.Cells.Item(4, 13).Value = -308.85715  # M4

$ws.Cells.Item(13, 8).Value = 8750  # H13
$ws.Cells.Item(13, 10).Value = 8750  # J13
$ws.Cells.Item(13, 12).Value = 8750  # L13
$ws.Cells.Item(13, 14).Value = -9038  # N13

$ws.Cells.Item(32, 8).Value = 3504.4546  # H32
$ws.Cells.Item(32, 9).Value = 2765.5134  # I32
$ws.Cells.Item(32, 10).Value = 7410.2856  # J32
$ws.Cells.Item(32, 11).Value = 2765.5134  # K32
$ws.Cells.Item(32, 12).Value = 7410.2856  # L32
$ws.Cells.Item(32, 13).Value = -2478.5134  # M32
$ws.Cells.Item(32, 14).Value = -7984.2856  # N32

$ws.Cells.Item(63, 8).Value = 6542.273  # H63
$ws.Cells.Item(63, 10).Value = 8867.666999999999  # J63
$ws.Cells.Item(63, 12).Value = 8867.666999999999  # L63
$ws.Cells.Item(63, 14).Value = -10239.667  # N63

$ws.Cells.Item(66, 8).Value = 6542.273  # H66
$ws.Cells.Item(66, 10).Value = 8867.666999999999  # J66
$ws.Cells.Item(66, 12).Value = 44338.335  # L66
$ws.Cells.Item(66, 14).Value = -51202.335  # N66

$ws.Cells.Item(131, 8).Value = 50000  # H131
$ws.Cells.Item(131, 10).Value = 50000  # J131
$ws.Cells.Item(131, 12).Value = 50000  # L131
$ws.Cells.Item(131, 14).Value = -60080  # N131 (new cell)

$ws.Cells.Item(132, 8).Value = 2605.8171  # H132
$ws.Cells.Item(132, 9).Value = 864.8840300000001  # I132
$ws.Cells.Item(132, 11).Value = 2594.65209  # K132
$ws.Cells.Item(132, 13).Value = -64.65209000000004  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3362.4375  # H20
$ws.Cells.Item(20, 9).Value = 2451.6667  # I20
$ws.Cells.Item(20, 10).Value = 3908.9  # J20
$ws.Cells.Item(20, 11).Value = 2451.6667  # K20
$ws.Cells.Item(20, 12).Value = 3908.9  # L20
$ws.Cells.Item(20, 13).Value = -2204.6667  # M20
$ws.Cells.Item(20, 14).Value = -4402.9  # N20

$ws.Cells.Item(22, 8).Value = 1  # H22
$ws.Cells.Item(22, 9).Value = 1  # I22
$ws.Cells.Item(22, 11).Value = 1  # K22
$ws.Cells.Item(22, 13).Value = 172  # M22

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 913093.6  # H31
$ws.Cells.Item(31, 9).Value = 1113336.6  # I31
$ws.Cells.Item(31, 10).Value = 12000  # J31
$ws.Cells.Item(31, 11).Value = 1113336.6  # K31
$ws.Cells.Item(31, 12).Value = 12000  # L31
$ws.Cells.Item(31, 13).Value = -1113041.6  # M31
$ws.Cells.Item(31, 14).Value = -12590  # N31

$ws.Cells.Item(34, 8).Value = 913093.6  # H34
$ws.Cells.Item(34, 9).Value = 1113336.6  # I34
$ws.Cells.Item(34, 10).Value = 12000  # J34
$ws.Cells.Item(34, 11).Value = 1113336.6  # K34
$ws.Cells.Item(34, 12).Value = 12000  # L34
$ws.Cells.Item(34, 13).Value = -1113134.6  # M34
$ws.Cells.Item(34, 14).Value = -12404  # N34

$ws.Cells.Item(111, 8).Value = 69249  # H111
$ws.Cells.Item(111, 10).Value = 69249  # J111
$ws.Cells.Item(111, 12).Value = 69249  # L111
$ws.Cells.Item(111, 14).Value = -77429  # N111

$ws.Cells.Item(141, 8).Value = 31086.111  # H141
$ws.Cells.Item(141, 10).Value = 101124.625  # J141
$ws.Cells.Item(141, 12).Value = 101124.625  # L141
$ws.Cells.Item(141, 14).Value = -111484.625  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 349.66666  # H3
$ws.Cells.Item(3, 9).Value = 349.66666  # I3
$ws.Cells.Item(3, 11).Value = 1048.99998  # K3
$ws.Cells.Item(3, 13).Value = -936.9999800000001  # M3

$ws.Cells.Item(13, 8).Value = 1621.25  # H13
$ws.Cells.Item(13, 9).Value = 939.7778  # I13
$ws.Cells.Item(13, 11).Value = 2819.3334  # K13
$ws.Cells.Item(13, 13).Value = -2651.3334  # M13

$ws.Cells.Item(34, 8).Value = 56058.047  # H34
$ws.Cells.Item(34, 9).Value = 174.14285  # I34
$ws.Cells.Item(34, 10).Value = 84000  # J34
$ws.Cells.Item(34, 11).Value = 522.4285500000001  # K34
$ws.Cells.Item(34, 12).Value = 252000  # L34
$ws.Cells.Item(34, 13).Value = -438.4285500000001  # M34
$ws.Cells.Item(34, 14).Value = -252168  # N34

$ws.Cells.Item(36, 8).Value = 1786.4286  # H36
$ws.Cells.Item(36, 9).Value = 1233  # I36
$ws.Cells.Item(36, 10).Value = 2201.5  # J36
$ws.Cells.Item(36, 11).Value = 3699  # K36
$ws.Cells.Item(36, 12).Value = 6604.5  # L36
$ws.Cells.Item(36, 13).Value = -3530  # M36
$ws.Cells.Item(36, 14).Value = -6942.5  # N36

$ws.Cells.Item(38, 8).Value = 31.428572  # H38
$ws.Cells.Item(38, 9).Value = 20.2  # I38
$ws.Cells.Item(38, 11).Value = 60.59999999999999  # K38
$ws.Cells.Item(38, 13).Value = 286.4  # M38

$ws.Cells.Item(50, 8).Value = 14904.143  # H50
$ws.Cells.Item(50, 9).Value = 76.333336  # I50
$ws.Cells.Item(50, 11).Value = 229.000008  # K50
$ws.Cells.Item(50, 13).Value = 251.999992  # M50

$ws.Cells.Item(53, 8).Value = 14904.143  # H53
$ws.Cells.Item(53, 9).Value = 76.333336  # I53
$ws.Cells.Item(53, 11).Value = 229.000008  # K53
$ws.Cells.Item(53, 13).Value = 251.999992  # M53

$ws.Cells.Item(64, 9).Value = 333334000  # I64
$ws.Cells.Item(64, 10).Value = 0  # J64
$ws.Cells.Item(64, 11).Value = 1000002000  # K64
$ws.Cells.Item(64, 12).Value = 0  # L64
$ws.Cells.Item(64, 13).Value = -1000001730  # M64
$ws.Cells.Item(64, 14).ClearContents()  # N64 (cell removed)

$ws.Cells.Item(67, 9).Value = 333334000  # I67
$ws.Cells.Item(67, 10).Value = 0  # J67
$ws.Cells.Item(67, 11).Value = 1000002000  # K67
$ws.Cells.Item(67, 12).Value = 0  # L67
$ws.Cells.Item(67, 13).Value = -1000001064  # M67
$ws.Cells.Item(67, 14).ClearContents()  # N67 (cell removed)

$ws.Cells.Item(92, 8).Value = 866.44446  # H92
$ws.Cells.Item(92, 9).Value = 523  # I92
$ws.Cells.Item(92, 10).Value = 964.5714  # J92
$ws.Cells.Item(92, 11).Value = 1569  # K92
$ws.Cells.Item(92, 12).Value = 2893.7142  # L92
$ws.Cells.Item(92, 13).Value = -321  # M92
$ws.Cells.Item(92, 14).Value = -5389.7142  # N92

$ws.Cells.Item(114, 8).Value = 1205.7059  # H114
$ws.Cells.Item(114, 10).Value = 1806.5555  # J114
$ws.Cells.Item(114, 12).Value = 5419.666499999999  # L114
$ws.Cells.Item(114, 14).Value = -11927.6665  # N114

$ws.Cells.Item(117, 8).Value = 671.0769  # H117
$ws.Cells.Item(117, 9).Value = 447.7143  # I117
$ws.Cells.Item(117, 10).Value = 931.6667  # J117
$ws.Cells.Item(117, 11).Value = 1343.1429  # K117
$ws.Cells.Item(117, 12).Value = 2795.0001  # L117
$ws.Cells.Item(117, 13).Value = 2098.8571  # M117
$ws.Cells.Item(117, 14).Value = -9679.000100000001  # N117

$ws.Cells.Item(122, 8).Value = 91907.55  # H122
$ws.Cells.Item(122, 10).Value = 100998.3  # J122
$ws.Cells.Item(122, 12).Value = 908984.7000000001  # L122
$ws.Cells.Item(122, 14).Value = -913884.7000000001  # N122

$ws.Cells.Item(131, 8).Value = 4412.2  # H131
$ws.Cells.Item(131, 9).Value = 958.1667  # I131
$ws.Cells.Item(131, 10).Value = 5892.5  # J131
$ws.Cells.Item(131, 11).Value = 2874.5001  # K131
$ws.Cells.Item(131, 12).Value = 17677.5  # L131
$ws.Cells.Item(131, 13).Value = 2165.4999  # M131
$ws.Cells.Item(131, 14).Value = -27757.5  # N131

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2924.4375  # H132
$ws.Cells.Item(132, 9).Value = 2423.8206  # I132
$ws.Cells.Item(132, 10).Value = 5093.778  # J132
$ws.Cells.Item(132, 11).Value = 7271.4618  # K132
$ws.Cells.Item(132, 12).Value = 15281.334  # L132
$ws.Cells.Item(132, 13).Value = -4741.4618  # M132
$ws.Cells.Item(132, 14).Value = -20341.334  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 34486800  # H122
$ws.Cells.Item(122, 10).Value = 6248.25  # J122
$ws.Cells.Item(122, 12).Value = 18744.75  # L122
$ws.Cells.Item(122, 14).Value = -23644.75  # N122

$ws.Cells.Item(132, 8).Value = 28622.184  # H132
$ws.Cells.Item(132, 9).Value = 1232.5172  # I132
$ws.Cells.Item(132, 11).Value = 3697.5516  # K132
$ws.Cells.Item(132, 13).Value = -1167.5516  # M132
